$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'CLfwZ777'
$ws.Range("B2").Value = 23080256
$ws.Range("C2").Value = 'kaqhund46'
$ws.Range("D2").Value = 'wY8$X%2n'
$ws.Range("F2").Value = 'SIJjvdpw'
$ws.Range("G2").Value = 'YBKM'
